# Update Pglyrp1-Trem1 LR-pair stats after recount of ligand/receptor-expressing
# cells (Dr Hou advice): Ligand-expressing cells & Receptor-expressing cells go
# from 1 to 3 for every row, and all derived expression/specificity metrics are
# recomputed accordingly across the whole data range (rows 2-11).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.995483333333333
$ws.Range("H2").Value = 5.98645
$ws.Range("I2").Value = 0.3457527559234298
$ws.Range("J2").Value = 0.3457527559234297
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.217721
$ws.Range("N2").Value = 9.653162999999999
$ws.Range("O2").Value = 0.2700934312193076
$ws.Range("P2").Value = 0.2700934312193076
$ws.Range("Q2").Value = 6.420908626816667
$ws.Range("R2").Value = 57.78817764135
$ws.Range("S2").Value = 0.09338554820089093
$ws.Range("T2").Value = 0.09338554820089091

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.995483333333333
$ws.Range("H3").Value = 5.98645
$ws.Range("I3").Value = 0.3457527559234298
$ws.Range("J3").Value = 0.3457527559234297
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.695641666666667
$ws.Range("N3").Value = 26.086925
$ws.Range("O3").Value = 0.7299065687806925
$ws.Range("P3").Value = 0.7299065687806925
$ws.Range("Q3").Value = 17.35200801847222
$ws.Range("R3").Value = 156.16807216625
$ws.Range("S3").Value = 0.2523672077225389
$ws.Range("T3").Value = 0.2523672077225388

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.073876333333333
$ws.Range("H4").Value = 3.221629
$ws.Range("I4").Value = 0.1860680545753899
$ws.Range("J4").Value = 0.1860680545753899
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.217721
$ws.Range("N4").Value = 9.653162999999999
$ws.Range("O4").Value = 0.2700934312193076
$ws.Range("P4").Value = 0.2700934312193076
$ws.Range("Q4").Value = 3.455434429169666
$ws.Range("R4").Value = 31.098909862527
$ws.Range("S4").Value = 0.05025575930056845
$ws.Range("T4").Value = 0.05025575930056845

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.073876333333333
$ws.Range("H5").Value = 3.221629
$ws.Range("I5").Value = 0.1860680545753899
$ws.Range("J5").Value = 0.1860680545753899
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 8.695641666666667
$ws.Range("N5").Value = 26.086925
$ws.Range("O5").Value = 0.7299065687806925
$ws.Range("P5").Value = 0.7299065687806925
$ws.Range("Q5").Value = 9.338043788980556
$ws.Range("R5").Value = 84.04239410082501
$ws.Range("S5").Value = 0.1358122952748215
$ws.Range("T5").Value = 0.1358122952748215

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.9363446666666667
$ws.Range("H6").Value = 2.809034
$ws.Range("I6").Value = 0.1622382625734142
$ws.Range("J6").Value = 0.1622382625734142
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.217721
$ws.Range("N6").Value = 9.653162999999999
$ws.Range("O6").Value = 0.2700934312193076
$ws.Range("P6").Value = 0.2700934312193076
$ws.Range("Q6").Value = 3.012895897171333
$ws.Range("R6").Value = 27.116063074542
$ws.Range("S6").Value = 0.04381948901351242
$ws.Range("T6").Value = 0.04381948901351242

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.9363446666666667
$ws.Range("H7").Value = 2.809034
$ws.Range("I7").Value = 0.1622382625734142
$ws.Range("J7").Value = 0.1622382625734142
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.695641666666667
$ws.Range("N7").Value = 26.086925
$ws.Range("O7").Value = 0.7299065687806925
$ws.Range("P7").Value = 0.7299065687806925
$ws.Range("Q7").Value = 8.142117697827778
$ws.Range("R7").Value = 73.27905928045
$ws.Range("S7").Value = 0.1184187735599018
$ws.Range("T7").Value = 0.1184187735599018

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.504935
$ws.Range("H8").Value = 4.514805
$ws.Range("I8").Value = 0.2607565871604841
$ws.Range("J8").Value = 0.2607565871604841
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.217721
$ws.Range("N8").Value = 9.653162999999999
$ws.Range("O8").Value = 0.2700934312193076
$ws.Range("P8").Value = 0.2700934312193076
$ws.Range("Q8").Value = 4.842460953134999
$ws.Range("R8").Value = 43.582148578215
$ws.Range("S8").Value = 0.07042864133921159
$ws.Range("T8").Value = 0.07042864133921159

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.504935
$ws.Range("H9").Value = 4.514805
$ws.Range("I9").Value = 0.2607565871604841
$ws.Range("J9").Value = 0.2607565871604841
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.695641666666667
$ws.Range("N9").Value = 26.086925
$ws.Range("O9").Value = 0.7299065687806925
$ws.Range("P9").Value = 0.7299065687806925
$ws.Range("Q9").Value = 13.086375491625
$ws.Range("R9").Value = 117.777379424625
$ws.Range("S9").Value = 0.1903279458212725
$ws.Range("T9").Value = 0.1903279458212725

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.2607776666666666
$ws.Range("H10").Value = 0.7823329999999999
$ws.Range("I10").Value = 0.04518433976728187
$ws.Range("J10").Value = 0.04518433976728187
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.217721
$ws.Range("N10").Value = 9.653162999999999
$ws.Range("O10").Value = 0.2700934312193076
$ws.Range("P10").Value = 0.2700934312193076
$ws.Range("Q10").Value = 0.8391097743643331
$ws.Range("R10").Value = 7.551987969278999
$ws.Range("S10").Value = 0.01220399336512417
$ws.Range("T10").Value = 0.01220399336512417

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.2607776666666666
$ws.Range("H11").Value = 0.7823329999999999
$ws.Range("I11").Value = 0.04518433976728187
$ws.Range("J11").Value = 0.04518433976728187
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 8.695641666666667
$ws.Range("N11").Value = 26.086925
$ws.Range("O11").Value = 0.7299065687806925
$ws.Range("P11").Value = 0.7299065687806925
$ws.Range("Q11").Value = 2.267629144002778
$ws.Range("R11").Value = 20.408662296025
$ws.Range("S11").Value = 0.0329803464021577
$ws.Range("T11").Value = 0.0329803464021577
